$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot in columns D (Price) and
# E (Volume(1h)) for rows 2-51. Some "price" strings look numeric (e.g.
# "0.3305"); Range("...").NumberFormat = "@" is applied first on those cells
# so Excel keeps them as literal text (preserving exact digits/trailing
# zeros) instead of silently re-parsing them as numbers.
$ws.Cells.Item(2, 4).Value = "30.331.87"
$ws.Cells.Item(2, 5).Value = "  -2.54%  "
$ws.Cells.Item(3, 4).Value = "1.934.61"
$ws.Cells.Item(3, 5).Value = "  -2.35%  "
$ws.Cells.Item(4, 5).Value = "  -0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "250.66"
$ws.Cells.Item(5, 5).Value = "  -1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.7097"
$ws.Cells.Item(6, 5).Value = "  -3.43%  "
$ws.Cells.Item(7, 5).Value = "  -0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3305"
$ws.Cells.Item(8, 5).Value = "  -2.17%  "
$ws.Cells.Item(9, 5).Value = "  +1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.07301"
$ws.Cells.Item(10, 5).Value = "  +2.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.8056"
$ws.Cells.Item(11, 5).Value = "  -2.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08068"
$ws.Cells.Item(12, 5).Value = "  -0.45%  "
$ws.Cells.Item(13, 4).Value = "1.935.38"
$ws.Cells.Item(13, 5).Value = "  -2.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.485"
$ws.Cells.Item(14, 5).Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "94.64"
$ws.Cells.Item(15, 5).Value = "  -4.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "15.13"
$ws.Cells.Item(16, 5).Value = "  -1.35%  "
$ws.Cells.Item(17, 4).Value = "30.333.81"
$ws.Cells.Item(17, 5).Value = "  -2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "253.58"
$ws.Cells.Item(18, 5).Value = "  -5.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000008214"
$ws.Cells.Item(19, 5).Value = "  -0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "5.795"
$ws.Cells.Item(20, 5).Value = "  -3.88%  "
$ws.Cells.Item(21, 4).Value = "2.190.02"
$ws.Cells.Item(21, 5).Value = "  -3.01%  "
$ws.Cells.Item(22, 5).Value = "  -0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "1.000"
$ws.Cells.Item(23, 5).Value = "  -0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.983"
$ws.Cells.Item(24, 5).Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.756"
$ws.Cells.Item(25, 5).Value = "  -1.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "164.98"
$ws.Cells.Item(26, 5).Value = "  +1.80%  "
$ws.Cells.Item(27, 5).Value = "  -1.58%  "
$ws.Cells.Item(28, 5).Value = "  -0.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.1286"
$ws.Cells.Item(29, 5).Value = "  -2.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.350"
$ws.Cells.Item(30, 5).Value = "  -2.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.540"
$ws.Cells.Item(31, 5).Value = "  -3.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.427"
$ws.Cells.Item(32, 5).Value = "  -3.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.175"
$ws.Cells.Item(33, 5).Value = "  -4.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.05206"
$ws.Cells.Item(34, 5).Value = "  -1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.264"
$ws.Cells.Item(35, 5).Value = "  -1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.7485"
$ws.Cells.Item(36, 5).Value = "  -3.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.793"
$ws.Cells.Item(37, 5).Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01969"
$ws.Cells.Item(38, 5).Value = "  -1.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.812"
$ws.Cells.Item(39, 5).Value = "  -2.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "78.94"
$ws.Cells.Item(40, 5).Value = "  -5.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "6.420"
$ws.Cells.Item(41, 5).Value = "  -5.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.4534"
$ws.Cells.Item(42, 5).Value = "  -1.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.021"
$ws.Cells.Item(43, 5).Value = "  -4.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.8447"
$ws.Cells.Item(44, 5).Value = "  -1.14%  "
$ws.Cells.Item(45, 5).Value = "  -0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "101.71"
$ws.Cells.Item(46, 5).Value = "  -2.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.809"
$ws.Cells.Item(47, 5).Value = "  -2.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "7.458"
$ws.Cells.Item(48, 5).Value = "  -2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "36.79"
$ws.Cells.Item(49, 5).Value = "  -1.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.4179"
$ws.Cells.Item(50, 5).Value = "  -2.65%  "
$ws.Cells.Item(51, 5).Value = "  -0.10%  "
